$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 181, shifting existing rows 181..236 down to 182..237
$ws.Rows.Item(181).Insert()

# Populate the newly inserted row 181 with the new record
$ws.Range("A181").Value = 4
$ws.Range("B181").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C181").Value = "Los Lagos"
$ws.Range("D181").Value = 44627
$ws.Range("E181").Value = 10
$ws.Range("F181").Value = 100112043
$ws.Range("G181").Value = "Pepino ensalada"
$ws.Range("H181").Value = "Sin especificar"
$ws.Range("I181").Value = "Primera"
$ws.Range("J181").Value = 100
$ws.Range("K181").Value = 21000
$ws.Range("L181").Value = 21000
$ws.Range("M181").Value = 21000
$ws.Range("N181").Value = "`$/caja 60 unidades"
$ws.Range("O181").Value = "Región de Arica y Parinacota"
$ws.Range("P181").Value = 350
$ws.Range("Q181").Value = 60
$ws.Range("R181").Value = "Hortaliza"
